$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sheet ---
$ws.Name = "HomeLandingTopicCards"

# --- Swap header columns D1/E1 (content + style) ---
# Before: D1 = LinkSelector (style w/o text-numfmt), E1 = TitleSelector (style w/ text-numfmt)
# After:  D1 = TitleSelector (style w/ text-numfmt), E1 = LinkSelector (style w/o text-numfmt)
$ws.Range("D1").Copy($ws.Range("H1"))
$ws.Range("E1").Copy($ws.Range("D1"))
$ws.Range("H1").Copy($ws.Range("E1"))
$ws.Range("H1").Clear()

# --- New data rows (3-13) ---
$rows = @(
    @(3,  "/",                              "Home",    "Guide",            ".guide-card .card h2",                    ".guide-card .card h2  + ul li a", 0, $false),
    @(4,  "/",                              "Home",    "Multimedia",       ".multimedia div[class*=feature-card] h3", ".multimedia div[class*=feature-card] h3", 0, $false),
    @(5,  "/",                              "Home",    "Thumbnail",        ".card-thumbnail h3 a",                    ".card-thumbnail h3 a", 0, $true),
    @(6,  "/espanol",                       "Home",    "Guide",            ".guide-card .card h2",                    ".guide-card .card h2  + ul li a", 0, $false),
    @(7,  "/espanol",                       "Home",    "Thumbnail",        ".card-thumbnail h3 a",                    ".card-thumbnail h3 a", 1, $true),
    @(8,  "/about-cancer",                  "Landing", "Feature",          ".feature-primary .feature-card h3",       ".feature-primary .feature-card h3", 0, $false),
    @(9,  "/about-cancer",                  "Landing", "SecondaryFeature", ".feature-secondary .feature-card h3",     ".feature-secondary .feature-card h3", 1, $true),
    @(10, "/espanol/cancer",                "Landing", "SecondaryFeature", ".feature-secondary .feature-card h3",     ".feature-secondary .feature-card h3", 2, $true),
    @(11, "/about-cancer/treatment",        "Topic",   "InlineCard",       "#cgvBody .feature-card a h3",             "#cgvBody .feature-card a h3", 0, $true),
    @(12, "/about-cancer/treatment",        "Topic",   "Thumbnail",        ".card-thumbnail h3 a",                    ".card-thumbnail h3 a", 2, $true),
    @(13, "/about-nci/organization/crchd",  "Topic",   "SlottedTopicCard", ".topic-feature .feature-card a h3",       ".topic-feature .feature-card a h3", 1, $true)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    if ($row[7]) {
        $ws.Range("D" + $r + ":E" + $r).NumberFormat = "@"
    } else {
        # Column E still carries a legacy column-level text format (inherited
        # from before the D/E swap); freshly-created cells in that column
        # pick it up automatically, so force them back to the default style.
        $ws.Cells.Item($r, 5).Style = "Normal"
    }
}

# --- Column widths (chars match the stored OOXML widths as closely as this engine's rounding allows) ---
$ws.Columns.Item(3).ColumnWidth = 16.5
$ws.Columns.Item(4).ColumnWidth = 36.16666666666667
$ws.Columns.Item(5).ColumnWidth = 36.16666666666667

# --- Selection moves to the next empty row ---
$ws.Range("A14").Select()

Write-Output "done"
